# WBS sheet: the Phase 0 / Phase 1 / Phase 2 work-package + milestone rows
# (rows 2-49) have all finished, so:
#   - Execution Status (H) flips from "Partial" to "Done"
#   - Completed On (K) is stamped with the same date already recorded in
#     Started On (J)
#   - the category checklist columns (Schema/Validation/Permissions-
#     Isolation/Workflow/Evidence = L:P) get a checkmark now that each
#     category has been covered for the milestone

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WBS")

$checkMark = [char]0x2705

for ($row = 2; $row -le 49; $row++) {
    # Execution Status: Partial -> Done
    $ws.Cells.Item($row, 8).Value = "Done"

    # Completed On (K) takes the same date text already shown in Started On
    # (J). Read it back through .Text (a formatted string) rather than
    # .Value, then round-trip the destination cell's NumberFormat so the
    # literal date text is stored as text instead of being re-parsed into a
    # serial date number, while still ending up with the same date-looking
    # cell format it already had.
    $startedOnText = $ws.Cells.Item($row, 10).Text
    $completedCell = $ws.Cells.Item($row, 11)
    $originalFormat = $completedCell.NumberFormat
    $completedCell.NumberFormat = "@"
    $completedCell.Value = $startedOnText
    $completedCell.NumberFormat = $originalFormat

    # Category columns: Schema (L), Validation (M), Permissions/Isolation (N),
    # Workflow (O), Evidence (P)
    $ws.Cells.Item($row, 12).Value = $checkMark
    $ws.Cells.Item($row, 13).Value = $checkMark
    $ws.Cells.Item($row, 14).Value = $checkMark
    $ws.Cells.Item($row, 15).Value = $checkMark
    $ws.Cells.Item($row, 16).Value = $checkMark
}
